$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.085.56'
$ws.Range("E2").Value = '  -3.74%  '

$ws.Range("D3").Value = '2.452.09'
$ws.Range("E3").Value = '  -3.02%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.84'
$ws.Range("E5").Value = '  +0.19%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '92.38'
$ws.Range("E6").Value = '  -8.41%  '

$ws.Range("E7").Value = '  -2.99%  '

$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("E9").Value = '  -5.80%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.97'
$ws.Range("E10").Value = '  -7.70%  '

$ws.Range("E11").Value = '  -4.05%  '

$ws.Range("E12").Value = '  -0.85%  '

$ws.Range("E13").Value = '  -5.62%  '

$ws.Range("D14").Value = '2.831.86'
$ws.Range("E14").Value = '  -2.90%  '

$ws.Range("D15").Value = '2.437.82'
$ws.Range("E15").Value = '  -6.36%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.59'
$ws.Range("E16").Value = '  -4.86%  '

$ws.Range("E17").Value = '  -4.08%  '

$ws.Range("D18").Value = '41.076.09'
$ws.Range("E18").Value = '  -3.73%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.24'
$ws.Range("E19").Value = '  -7.20%  '

$ws.Range("E20").Value = '  -4.28%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.10'
$ws.Range("E21").Value = '  -9.12%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.70'
$ws.Range("E22").Value = '  -2.87%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '234.56'
$ws.Range("E23").Value = '  -3.60%  '

$ws.Range("E24").Value = '  -4.98%  '

$ws.Range("E25").Value = '  +0.34%  '

$ws.Range("E26").Value = '  -7.05%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.71'
$ws.Range("E27").Value = '  -6.95%  '

$ws.Range("E28").Value = '  -5.74%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.52'
$ws.Range("E29").Value = '  -5.78%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.41'
$ws.Range("E30").Value = '  -8.17%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '150.41'
$ws.Range("E31").Value = '  -4.45%  '

$ws.Range("E32").Value = '  -5.50%  '

$ws.Range("E33").Value = '  -5.17%  '

$ws.Range("E34").Value = '  -3.31%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0732'
$ws.Range("E35").Value = '  -6.53%  '

$ws.Range("E36").Value = '  -5.84%  '

$ws.Range("B37").Value = 'Celestia'
$ws.Range("C37").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '16.68'
$ws.Range("E37").Value = '  -6.90%  '

$ws.Range("B38").Value = 'ARBITRUM'
$ws.Range("C38").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.84'
$ws.Range("E38").Value = '  -6.78%  '

$ws.Range("E39").Value = '  -3.71%  '

$ws.Range("E40").Value = '  -8.83%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.12'
$ws.Range("E41").Value = '  -2.11%  '

$ws.Range("E42").Value = '  +0.18%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '19.93'

$ws.Range("D44").Value = '1.958.94'
$ws.Range("E44").Value = '  -2.35%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0280'
$ws.Range("E45").Value = '  -6.64%  '

$ws.Range("E46").Value = '  -9.73%  '

$ws.Range("E47").Value = '  -4.60%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '70.13'
$ws.Range("E48").Value = '  -2.71%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '95.87'
$ws.Range("E49").Value = '  -5.05%  '

$ws.Range("E50").Value = '  -7.60%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.31'
$ws.Range("E51").Value = '  -7.25%  '
